$wb = $excel.ActiveWorkbook

# Rename sheets (by index, to preserve order matching workbook.xml sheet order)
$wb.Worksheets.Item(1).Name = "GNG_TO-16504777734963791"
$wb.Worksheets.Item(2).Name = "NB_TO-1650477775785389"
$wb.Worksheets.Item(3).Name = "RS_TO-16504777757883766"
$wb.Worksheets.Item(4).Name = "TOL_TO-1650477775864383"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16504777759593773"

# Sheet 1 (GNG)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16504777734483786.csv"
$ws1.Range("B3").Value = "GNG_stims-16504777734633775.csv"
$ws1.Range("B4").Value = "go_stims-16504777734653876.csv"
$ws1.Range("B5").Value = "GNG_stims-16504777734963791.csv"

# Sheet 2 (NB)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16504777743413944.csv"
$ws2.Range("B3").Value = "ZB-match_7-165047777407838.csv"
$ws2.Range("B4").Value = "TB-16504777757673814.csv"
$ws2.Range("B5").Value = "TB-1650477775162422.csv"
$ws2.Range("B6").Value = "ZB-match_7-16504777741903777.csv"
$ws2.Range("B7").Value = "OB-16504777751003838.csv"
$ws2.Range("B8").Value = "ZB-match_3-1650477773629375.csv"
$ws2.Range("B9").Value = "OB-16504777745364.csv"
$ws2.Range("B10").Value = "TB-16504777752233841.csv"

# Sheet 3 (RS)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# Sheet 4 (TOL)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16504777758174133.csv"
$ws4.Range("B3").Value = "ZM_stims-16504777757923777.csv"
$ws4.Range("B4").Value = "MM_stims-16504777758484123.csv"
$ws4.Range("B5").Value = "ZM_stims-16504777758183753.csv"
$ws4.Range("B6").Value = "MM_stims-16504777758633842.csv"
$ws4.Range("B7").Value = "ZM_stims-16504777758503773.csv"

# Sheet 5 (vSAT)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16504777758703861.csv"
$ws5.Range("B3").Value = "SAT_stims-16504777758963819.csv"
$ws5.Range("B4").Value = "vSAT_stims-16504777759123855.csv"
$ws5.Range("B5").Value = "vSAT_stims-16504777759433832.csv"
